$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 210; all rows from 210 downward shift down by one
# (old row 210 becomes row 211, ..., old row 297 becomes row 298).
$ws.Rows.Item(210).Insert()

# Populate the newly inserted row 210 with the new data record.
$ws.Range("A210").Value = 9
$ws.Range("B210").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C210").Value = "Metropolitana"
$ws.Range("D210").Value = 44795
$ws.Range("E210").Value = 13
$ws.Range("F210").Value = 100112001
$ws.Range("G210").Value = "Berenjena"
$ws.Range("H210").Value = "Sin especificar"
$ws.Range("I210").Value = "Primera"
$ws.Range("J210").Value = 70
$ws.Range("K210").Value = 10000
$ws.Range("L210").Value = 12000
$ws.Range("M210").Value = 11000
$ws.Range("N210").Value = "$/caja 50 unidades"
$ws.Range("O210").Value = "Región de Arica y Parinacota"
$ws.Range("P210").Value = 220
$ws.Range("Q210").Value = 50
$ws.Range("R210").Value = "Hortaliza"
